# Commit: "added org update and delete test case"
#
# DriveHealthTestData.xlsx gains several new "AUTO_ORG_*" org-name fixtures
# in the shared-string pool (used elsewhere by the test suite for update /
# delete scenarios), and the sample row on the "Data" sheet is repointed at
# one of the freshly added names so the single-row fixture exercises the
# new org value.
#
# Net effect on the worksheet itself: cell A2 on the "Data" sheet changes
# from "AUTO_ORG_ERZYN" to "AUTO_ORG_USPHR".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AUTO_ORG_USPHR"
